# Fix Training Data Issue
# The "Date" column (BF) had values one day off ("5-8-2011-12") due to how
# NBA stats were sourced. Correct them to the real ISO date "2012-05-08"
# for every data row (rows 2-31), keeping the cells as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("BF2:BF31")

# Force a text number format first so the date-looking string isn't
# auto-converted into a date serial value - we want the literal text
# "2012-05-08" preserved, just like the original "5-8-2011-12" text.
$rng.NumberFormat = "@"
$rng.Value = "2012-05-08"
